# Auto-generated edit script: updates market-price derived columns (H:N)
# across multiple crafting-leve sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 891.8
$ws.Range("I6").Value = 164.75
$ws.Range("J6").Value = 3800
$ws.Range("K6").Value = 494.25
$ws.Range("L6").Value = 11400
$ws.Range("M6").Value = -382.25
$ws.Range("N6").Value = -11624
# Row 64
$ws.Range("H64").Value = 2831.125
$ws.Range("I64").Value = 2699.8572
$ws.Range("J64").Value = 2933.2222
$ws.Range("K64").Value = 2699.8572
$ws.Range("L64").Value = 2933.2222
$ws.Range("M64").Value = -2451.8572
$ws.Range("N64").Value = -3429.2222
# Row 67
$ws.Range("H67").Value = 2831.125
$ws.Range("I67").Value = 2699.8572
$ws.Range("J67").Value = 2933.2222
$ws.Range("K67").Value = 2699.8572
$ws.Range("L67").Value = 2933.2222
$ws.Range("M67").Value = -1841.8572
$ws.Range("N67").Value = -4649.2222
# Row 86
$ws.Range("H86").Value = 5750
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 7333.3335
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 7333.3335
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -9579.333500000001
# Row 89
$ws.Range("H89").Value = 5750
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 7333.3335
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 36666.6675
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -47898.6675
# Row 127
$ws.Range("H127").Value = 2152.4546
$ws.Range("I127").Value = 1197
$ws.Range("J127").Value = 2248
$ws.Range("K127").Value = 3591
$ws.Range("L127").Value = 6744
$ws.Range("M127").Value = 1369
$ws.Range("N127").Value = -16664
# Row 134
$ws.Range("H134").Value = 51062
$ws.Range("J134").Value = 51062
$ws.Range("L134").Value = 51062
$ws.Range("N134").Value = -61202
# Row 137
$ws.Range("H137").Value = 2941
$ws.Range("I137").Value = 2578.9092
$ws.Range("J137").Value = 4389.364
$ws.Range("K137").Value = 7736.7276
$ws.Range("L137").Value = 13168.092
$ws.Range("M137").Value = -5186.7276
$ws.Range("N137").Value = -18268.092

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2221.8
$ws.Range("I61").Value = 1731.0952
$ws.Range("J61").Value = 2957.8572
$ws.Range("K61").Value = 1731.0952
$ws.Range("L61").Value = 2957.8572
$ws.Range("M61").Value = -1519.0952
$ws.Range("N61").Value = -3381.8572
# Row 122
$ws.Range("H122").Value = 2300.4583
$ws.Range("I122").Value = 1122.8334
$ws.Range("J122").Value = 5833.3335
$ws.Range("K122").Value = 3368.5002
$ws.Range("L122").Value = 17500.0005
$ws.Range("M122").Value = -918.5001999999999
$ws.Range("N122").Value = -22400.0005
# Row 132
$ws.Range("H132").Value = 2730.7322
$ws.Range("I132").Value = 2000.5814
$ws.Range("J132").Value = 5145.846
$ws.Range("K132").Value = 6001.7442
$ws.Range("L132").Value = 15437.538
$ws.Range("M132").Value = -3471.7442
$ws.Range("N132").Value = -20497.538
# Row 133
$ws.Range("H133").Value = 40808.75
$ws.Range("J133").Value = 40808.75
$ws.Range("L133").Value = 40808.75
$ws.Range("N133").Value = -45868.75
# Row 136
$ws.Range("H136").Value = 2221.8
$ws.Range("I136").Value = 1731.0952
$ws.Range("J136").Value = 2957.8572
$ws.Range("K136").Value = 5193.2856
$ws.Range("L136").Value = 8873.571599999999
$ws.Range("M136").Value = -2643.2856
$ws.Range("N136").Value = -13973.5716
# Row 137
$ws.Range("H137").Value = 53750
$ws.Range("J137").Value = 53750
$ws.Range("L137").Value = 53750
$ws.Range("N137").Value = -63950
# Row 140
$ws.Range("H140").Value = 50076.668
$ws.Range("J140").Value = 50076.668
$ws.Range("L140").Value = 50076.668
$ws.Range("N140").Value = -60436.668

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1895.6897
$ws.Range("I86").Value = 1844.8182
$ws.Range("J86").Value = 2055.5715
$ws.Range("K86").Value = 1844.8182
$ws.Range("L86").Value = 2055.5715
$ws.Range("M86").Value = -721.8181999999999
$ws.Range("N86").Value = -4301.5715
# Row 89
$ws.Range("H89").Value = 1895.6897
$ws.Range("I89").Value = 1844.8182
$ws.Range("J89").Value = 2055.5715
$ws.Range("K89").Value = 9224.091
$ws.Range("L89").Value = 10277.8575
$ws.Range("M89").Value = -3608.091
$ws.Range("N89").Value = -21509.8575
# Row 94
$ws.Range("H94").Value = 1135
$ws.Range("I94").Value = 1140
$ws.Range("K94").Value = 1140
$ws.Range("M94").Value = -689
# Row 107
$ws.Range("H107").Value = 564.8276
$ws.Range("I107").Value = 574.8929000000001
$ws.Range("J107").Value = 283
$ws.Range("K107").Value = 574.8929000000001
$ws.Range("L107").Value = 283
$ws.Range("M107").Value = 1345.1071
$ws.Range("N107").Value = -4123
# Row 134
$ws.Range("H134").Value = 2075.8735
$ws.Range("I134").Value = 1195.5883
$ws.Range("J134").Value = 7517.636
$ws.Range("K134").Value = 3586.7649
$ws.Range("L134").Value = 22552.908
$ws.Range("M134").Value = -1051.7649
$ws.Range("N134").Value = -27622.908
# Row 137
$ws.Range("H137").Value = 32965.715
$ws.Range("J137").Value = 32965.715
$ws.Range("L137").Value = 32965.715
$ws.Range("N137").Value = -43165.715

$ws = $wb.Worksheets.Item("CRP")
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 444.6
$ws.Range("I81").Value = 444.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1333.8
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -210.8000000000002
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 444.6
$ws.Range("I84").Value = 444.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4001.4
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 1614.6
$ws.Range("N84").ClearContents()
# Row 113
$ws.Range("H113").Value = 558.7344000000001
$ws.Range("I113").Value = 556.9778
$ws.Range("J113").Value = 562.8946999999999
$ws.Range("K113").Value = 1670.9334
$ws.Range("L113").Value = 1688.6841
$ws.Range("M113").Value = 499.0666000000001
$ws.Range("N113").Value = -6028.6841
# Row 137
$ws.Range("H137").Value = 4672.5454
$ws.Range("I137").Value = 3305.8
$ws.Range("J137").Value = 5811.5
$ws.Range("K137").Value = 9917.400000000001
$ws.Range("L137").Value = 17434.5
$ws.Range("M137").Value = -4817.400000000001
$ws.Range("N137").Value = -27634.5

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2566.4546
$ws.Range("I132").Value = 1565.9166
$ws.Range("J132").Value = 5234.5557
$ws.Range("K132").Value = 4697.7498
$ws.Range("L132").Value = 15703.6671
$ws.Range("M132").Value = -2167.7498
$ws.Range("N132").Value = -20763.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5252.0527
$ws.Range("I40").Value = 5234.5713
$ws.Range("K40").Value = 5234.5713
$ws.Range("M40").Value = -5098.5713

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 33410018
$ws.Range("I62").Value = 71432060
$ws.Range("K62").Value = 71432060
$ws.Range("M62").Value = -71431436
# Row 65
$ws.Range("H65").Value = 33410018
$ws.Range("I65").Value = 71432060
$ws.Range("K65").Value = 357160300
$ws.Range("M65").Value = -357157180
# Row 81
$ws.Range("H81").Value = 29221852
$ws.Range("I81").Value = 29221852
$ws.Range("K81").Value = 58443704
$ws.Range("M81").Value = -58442643
# Row 84
$ws.Range("H84").Value = 29221852
$ws.Range("I84").Value = 29221852
$ws.Range("K84").Value = 292218520
$ws.Range("M84").Value = -292213216
# Row 130
$ws.Range("H130").Value = 39817
$ws.Range("J130").Value = 39817
$ws.Range("L130").Value = 39817
$ws.Range("N130").Value = -49857
# Row 132
$ws.Range("H132").Value = 5651164.5
$ws.Range("I132").Value = 484.12122
$ws.Range("J132").Value = 12823182
$ws.Range("K132").Value = 1452.36366
$ws.Range("L132").Value = 38469546
$ws.Range("M132").Value = 1077.63634
$ws.Range("N132").Value = -38474606
Write-Output "Updated currentAveragePrice/LevePrice/LeveProfit columns (H:N) for scheduled market-data refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
